# Weekly update: insert a new data row at the top of the data block (row 13),
# pushing the existing rows 13-36 down to 14-37.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 13 (shifts rows 13..36 down to 14..37)
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with this week's data
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Vega Monumental Concepción"
$ws.Range("C13").Value = "Bíobío"
$ws.Range("D13").Value = 44680
$ws.Range("E13").Value = 8
$ws.Range("F13").Value = 100112013
$ws.Range("G13").Value = "Alcachofa"
$ws.Range("H13").Value = "Española"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 50
$ws.Range("K13").Value = 18000
$ws.Range("L13").Value = 19000
$ws.Range("M13").Value = 18600
$ws.Range("N13").Value = "$/caja 30 unidades"
$ws.Range("O13").Value = "Provincia de Limarí"
$ws.Range("P13").Value = 620
$ws.Range("Q13").Value = 30
$ws.Range("R13").Value = "Hortaliza"
